$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- E5: fix the superscript annotation text from " $^d, e$" to " $^{d, e}$" ---
# Cell currently holds rich text: run1 = "0.01-1" (default font),
# run2 = " $^d, e$" (superscript Arial 11). We need to insert "{" right after
# the "^" and "}" right before the trailing "$", keeping the superscript
# formatting on the whole annotation run.
$cell = $ws.Range("E5")

# Insert from the end first so earlier character offsets stay valid.
$cell.Characters(14, 1).Text = "}$"   # turn trailing "$" into "}$"
$cell.Characters(9, 1).Text = "^{"    # turn "^" into "^{"

# Re-apply the superscript run formatting over the whole annotation
# (" $^{d, e}$", now 10 chars starting at position 7) so it keeps matching
# the rest of the superscript text in the table.
$run = $cell.Characters(7, 10)
$run.Font.Superscript = $true
$run.Font.Size = 11
$run.Font.Name = "Arial"

# --- Move the active selection from E3 to E5 ---
$ws.Range("E5").Select() | Out-Null
